$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting old D:K to E:L
$ws.Columns.Item(4).Insert()

# Copy formatting from the (shifted) old D column (now E) into new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest period data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 96400
$ws.Range("D9").Value = 41100
$ws.Range("D10").Value = 55400
$ws.Range("D12").Value = 55400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 22400
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 252100
$ws.Range("D18").Value = -155700
$ws.Range("D20").Value = 41700
$ws.Range("D21").Value = -97100
$ws.Range("D22").Value = 43200
$ws.Range("D23").Value = -157200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -157200
$ws.Range("D27").Value = -157200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -41700
$ws.Range("D33").Value = -157200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -157200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 81800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 22500
$ws.Range("D44").Value = 41300
$ws.Range("D45").Value = 3800
$ws.Range("D46").Value = 149500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1600
$ws.Range("D49").Value = 229200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 61300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 441600
$ws.Range("D57").Value = 16800
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 116800
$ws.Range("D60").Value = 133600
$ws.Range("D61").Value = 110500
$ws.Range("D62").Value = 7400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 251500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -719800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 190100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -157200
$ws.Range("D83").Value = 16900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -171500
$ws.Range("D91").Value = -1700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -170100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 295000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -46600

# A handful of historical cells were corrected as part of this update
$ws.Range("F21").Value = -116200
$ws.Range("F83").Value = 500
$ws.Range("F89").Value = -70600
$ws.Range("E91").Value = -800
$ws.Range("F91").Value = -500
$ws.Range("F94").Value = -500
$ws.Range("F100").Value = 52300
$ws.Range("F102").Value = -18700
